# Fixes GitHub issue #138 ("cannot convert Cell value to Nullable<T>")
# Normalizes the sample data so every data row shares the same date and
# product name (previously rows 4/7 used distinct dates/products that
# triggered the bug).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: date 2021-03-02 -> 2021-03-01, product "測試商品2" -> "測試商品1"
$ws.Range("A4").Value = 44256
$ws.Range("D4").Value = "測試商品1"

# Row 7: date 2021-03-03 -> 2021-03-01, product "測試商品3" -> "測試商品1"
$ws.Range("A7").Value = 44256
$ws.Range("D7").Value = "測試商品1"

# Move the active selection from A2 to B1, matching the saved workbook view.
$ws.Range("B1").Select()
